$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 31313
$ws.Range("H3").Value = 11

# Update the selected cell / active selection
$ws.Range("G6").Select()
